# "Included Multi User test data"
# Add a new data row (row 3) to the Request-For-Information transmittals
# test-data sheet, mirroring the existing single-user row (row 2) but with
# a "To" value that carries multiple users joined by "@@".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "AutoTestAdmin@@AutoTestUser"
$ws.Range("C3").Value = "New Transmittal from Automation"
$ws.Range("D3").Value = "UnTick"
$ws.Range("E3").Value = "Correspondence"
$ws.Range("F3").Value = "Request for Information"
$ws.Range("M3").Value = "Comments for Request for Information"
